# Update the cryptocurrency price (D) / 1h volume change (E) figures, and swap the
# WrappedEther/Polkadot row 12<->13 name, link, price and volume values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.428.23'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '''  +1.30%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = '''1.677.07'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  +2.38%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = '''216.70'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  +1.30%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = '''0.5311'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '''  +1.14%  '
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = '''  -0.03%  '
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = '''  +3.44%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = '''0.06394'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '''  +1.32%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = '''21.67'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '''  +4.63%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = '''0.07811'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '''  +1.88%  '
$ws.Range("E11").ClearFormats()
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '''4.508'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '''  +1.90%  '
$ws.Range("E12").ClearFormats()
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '''1.658.90'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  +1.21%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = '''0.5557'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '''  +0.84%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = '''0.0₅8309'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '''  +2.19%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = '''65.56'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '''  +0.57%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = '''26.473.81'
$ws.Range("D17").ClearFormats()
$ws.Range("E18").Value = '''  -0.06%  '
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = '''4.731'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '''  +0.67%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = '''193.40'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '''  +2.36%  '
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = '''10.29'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '''  +1.40%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = '''6.340'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '''  +2.82%  '
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = '''  +0.00%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = '''142.28'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '''  -2.75%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = '''0.1286'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '''  +5.56%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = '''7.396'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '''  -0.23%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = '''16.21'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '''  +2.24%  '
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = '''1.430'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '''  +1.46%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = '''0.06229'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '''  +3.77%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = '''1.273'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '''  +1.21%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = '''3.607'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '''  +4.68%  '
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = '''  +0.89%  '
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = '''  +2.11%  '
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = '''  +1.87%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = '''0.6136'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '''  +6.97%  '
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = '''  +1.32%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = '''2.781'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '''  +0.68%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = '''6.155'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '''  +8.13%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = '''0.01628'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '''  +0.61%  '
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = '''1.080.68'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '''  +3.79%  '
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = '''0.8636'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '''  +1.20%  '
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = '''  -0.13%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = '''100.18'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '''  -0.50%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = '''1.821.98'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '''  +1.97%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = '''57.06'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '''  +2.95%  '
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = '''8.134'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '''  +0.39%  '
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = '''  -0.21%  '
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = '''  -3.02%  '
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = '''0.05209'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '''  +0.70%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = '''1.470'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '''  +5.71%  '
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = '''6.012'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '''  +1.66%  '
$ws.Range("E51").ClearFormats()
